$wb = $excel.ActiveWorkbook

# --- Rename the "Value lists" sheet to "abc list" and rewrite its
#     contents from the X,Y,Z list to the A,B,C list. ---
$abcList = $wb.Worksheets.Item("Value lists")
$abcList.Range("A1").Value = "A"
$abcList.Range("A2").Value = "B"
$abcList.Range("A3").Value = "C"
$abcList.Name = "abc list"

# --- Add a brand-new "xyz list" sheet after the last existing sheet,
#     carrying the X,Y,Z values that used to live on "Value lists". ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$xyzList = $wb.Worksheets.Add($null, $lastSheet)
$xyzList.Name = "xyz list"
$xyzList.Range("A1").Value = "X"
$xyzList.Range("A2").Value = "Y"
$xyzList.Range("A3").Value = "Z"

# --- Point the data validations on the first sheet at literal lists
#     instead of the old 'Value lists' range reference. ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A2:A1048576").Validation.Modify(3, 1, 1, '"A,B,C"')
$ws1.Range("B2:B1048576").Validation.Modify(3, 1, 1, '"X,Y,Z"')
